# Refresh the Coinranking "cryptos" snapshot (Price / Volume(1h) columns)
# exactly as published by the "Updated cryptos list ... with GitHub Actions"
# commit. Every Price/Volume cell in this sheet is plain text (pre-formatted
# strings like "310.58" or "  -1.04%  "), not a numeric cell, so the refresh
# writes text back. Column D in particular holds values that LOOK like plain
# numbers ("310.58", "1.000", "5.310", ...); Excel's Range.Value setter
# auto-detects those as doubles and would silently normalise/trim them (e.g.
# "1.000" -> 1, "5.310" -> 5.31), so every D-column write is preceded by
# NumberFormat = '@' to pin the cell to Text and preserve the exact string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.890.31'
$ws.Range('E2').Value = '  -1.66%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.810.48'
$ws.Range('E3').Value = '  -0.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '310.58'
$ws.Range('E5').Value = '  -1.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4625'
$ws.Range('E7').Value = '  +3.46%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3754'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07425'
$ws.Range('E9').Value = '  -1.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8773'
$ws.Range('E10').Value = '  -0.93%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.44'
$ws.Range('E11').Value = '  -2.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.784.36'
$ws.Range('E12').Value = '  -2.43%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.354'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.545'
$ws.Range('E14').Value = '  -3.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07049'
$ws.Range('E15').Value = '  -0.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.41'
$ws.Range('E16').Value = '  -2.72%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008754'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.74'
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '26.886.84'
$ws.Range('E21').Value = '  -1.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.310'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.85'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.997.03'
$ws.Range('E24').Value = '  -2.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.923'
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.53'
$ws.Range('E26').Value = '  -0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.56'
$ws.Range('E27').Value = '  -0.02%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.151'
$ws.Range('E28').Value = '  -9.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.301'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.59'
$ws.Range('E30').Value = '  -1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08911'
$ws.Range('E31').Value = '  +1.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7706'
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.166'
$ws.Range('E33').Value = '  -2.57%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.481'
$ws.Range('E34').Value = '  -0.68%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.904'
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.119'
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01958'
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05239'
$ws.Range('E39').Value = '  -1.73%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.431'
$ws.Range('E40').Value = '  +4.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.281'
$ws.Range('E41').Value = '  -1.44%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5348'
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.910'
$ws.Range('E43').Value = '  +1.78%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1661'
$ws.Range('E44').Value = '  -3.74%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.570'
$ws.Range('E45').Value = '  -1.98%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5068'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.32'
$ws.Range('E47').Value = '  -3.69%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '104.33'
$ws.Range('E48').Value = '  -1.37%  '
$ws.Range('E49').Value = '  -0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.665'
$ws.Range('E50').Value = '  -2.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06331'
$ws.Range('E51').Value = '  -0.73%  '
